{"js": "// Insert a new paragraph right after the paragraph that ends with\n// \"You can find this POC in the POC folder.\" and before the trailing\n// empty paragraphs at the end of the document body.\n\nconst body = context.document.body;\nconst searchResults = body.search(\"You can find this POC in the POC folder.\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target sentence in the document body.\");\n}\n\n// The search hit lives on the target paragraph; get that paragraph and\n// insert the new one right after it.\nconst hitRange = searchResults.items[0];\nconst targetParagraph = hitRange.paragraphs.getFirst();\ntargetParagraph.insertParagraph(\n  \"After integrating this concept in my project I have all the data required to make the calculations and to send to the microcontroller.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$target = \"You can find this POC in the POC folder.\"\n$newText = \"After integrating this concept in my project I have all the data required to make the calculations and to send to the microcontroller.\"\n\n# Locate the paragraph that ends with the target sentence (the last\n# paragraph of real content, right before the trailing empty paragraphs).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*$target*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $para = $d.Paragraphs.Item($targetIndex)\n    # Insert a brand-new (empty) paragraph right after it, then fill it in.\n    $para.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($targetIndex + 1)\n    $newPara.Range.Text = $newText\n}\n"}
